# Apply textual replacements described by the diff.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-11-30 Sunday" "2025-12-01 Monday"

Replace-Text "815×6=" "573×4="
Replace-Text "412×9=" "628×3="
Replace-Text "471×6=" "177×3="
Replace-Text "543×6=" "919×6="
Replace-Text "303×5=" "564×5="

Replace-Text "577×6=" "507×2="
Replace-Text "131×3=" "873×5="
Replace-Text "460×7=" "799×5="
Replace-Text "126×4=" "221×3="
Replace-Text "783×8=" "548×8="

Replace-Text "706×4=" "630×2="
Replace-Text "458×5=" "758×7="
Replace-Text "673×5=" "427×3="
Replace-Text "324×6=" "516×7="
Replace-Text "586×9=" "746×3="

Replace-Text "334×5=" "815×4="
Replace-Text "638×5=" "889×7="
Replace-Text "228×7=" "179×4="
Replace-Text "418×8=" "916×5="
Replace-Text "279×5=" "433×7="

Replace-Text "849×7=" "654×6="
Replace-Text "937×6=" "349×3="
Replace-Text "651×7=" "323×6="
Replace-Text "679×9=" "536×2="
Replace-Text "773×5=" "395×8="
